$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: PRODUTO - update description text (plain text, no numeric coercion risk)
$ws.Range("B2").Value = "CANTONEIRA PLASTICA 20x20x42"

# Row 3: QTDE - update quantity. Use a text formula + paste-as-values round trip
# so the numeric-looking string "150" stays stored as text (matching the
# original inline-string typing) instead of being coerced to a number, and
# without disturbing the cell's existing style/number format.
$ws.Range("B3").Formula = '="150"'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)

# Bump the QTDE font size - it now shares the larger font used by CODIGO
$ws.Range("B3").Font.Size = 150

# Row 4: CODIGO - update code value the same way, and keep its large font
$ws.Range("B4").Formula = '="1018"'
$ws.Range("B4").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$ws.Range("B4").Font.Size = 150

# Row 5: VOL - value unchanged ("10"), nothing to do

$excel.CutCopyMode = $false
